# AFDP-3349: Add new access control rules for DocumentRepository
# - Rename existing "Document Repository-default read access" rule text
#   (adds a space before the dash) on row 59.
# - Add 8 new DocumentRepository access-control rule rows (60-67) mirroring
#   the layout already used for Complaint / Case File / Task / Folder rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# 1. Fix up the existing "Document Repository" default read access rule
#    on row 59: the rule name gains a space before the dash.
# ---------------------------------------------------------------------
$ws.Range("B59").Value = "Document Repository -default read access"

# ---------------------------------------------------------------------
# 2. Helper data describing the new rows (60-67).
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 60; Height = 45; B = "DocumentRepository – Anybody can add comments";
       C = "DOC_REPO"; D = $null; G = "grant addComment to *";
       ColStyles = @{ B = "16"; C = "16"; D = "16"; E = "16"; F = "16"; G = "16" } },

    @{ Row = 61; Height = 30; B = "DocumentRepository – Lockout No Access Users";
       C = "DOC_REPO"; D = $null; G = "mandatory deny read to No Access";
       ColStyles = @{ B = "16"; C = "16"; D = "16"; E = "16"; F = "16"; G = "16" } },

    @{ Row = 62; Height = 30; B = "DocumentRepository – Anybody can add tag";
       C = "DOC_REPO"; D = $null; G = "grant addTag to *";
       ColStyles = @{ B = "16"; C = "16"; D = "16"; E = "16"; F = "16"; G = "16" } },

    @{ Row = 63; Height = 30; B = "DocumentRepository – Anybody can subscribe";
       C = "DOC_REPO"; D = $null; G = "grant subscribe to *";
       ColStyles = @{ B = "16"; C = "16"; D = "16"; E = "16"; F = "16"; G = "16" } },

    @{ Row = 64; Height = 30; B = "DocumentRepository – Restricted Flag";
       C = "DOC_REPO"; D = "restricted"; G = "deny read to *";
       ColStyles = @{ B = "16"; C = "18"; D = "18"; E = "18"; F = "18"; G = "18" } },

    @{ Row = 65; Height = 45; B = "DocumentRepository – Only participants can add files";
       C = "DOC_REPO"; D = $null; G = "grant add file to assignee, co-owner, supervisor, owning group, approver, collaborator";
       ColStyles = @{ B = "16"; C = "18"; D = "17"; E = "18"; F = "18"; G = "16" } },

    @{ Row = 66; Height = 45; B = "DocumentRepository – Only participants can save";
       C = "DOC_REPO"; D = $null; G = "grant save to assignee, co-owner, supervisor, owning group, approver, collaborator";
       ColStyles = @{ B = "16"; C = "18"; D = "17"; E = "18"; F = "18"; G = "16" } },

    @{ Row = 67; Height = 60; B = "DocumentRepository – Only participants can upload or replace files";
       C = "DOC_REPO"; D = $null; G = "grant uploadOrReplaceFile to assignee, co-owner, supervisor, owning group, approver, collaborator, reader";
       ColStyles = @{ B = "16"; C = "18"; D = "20"; E = "18"; F = "18"; G = "16" } }
)

# Reference cells (existing, already-styled) used as formatting sources so
# that the new cells reuse the same cellXf entries instead of creating new
# (duplicate) styles.
$styleSource = @{
    "16" = "B59"
    "18" = "C59"
    "17" = "D48"
    "20" = "D52"
}

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy cell formatting column by column so each cell ends up with the
    # same style index as the equivalent cell in the existing table.
    foreach ($col in @("B", "C", "D", "E", "F", "G")) {
        $styleId = $r.ColStyles[$col]
        $srcCell = $styleSource[$styleId]
        $ws.Range($srcCell).Copy()
        $ws.Range("$col$rowNum").PasteSpecial($xlPasteFormats)
    }

    # Now fill in the values.
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    if ($r.D) {
        $ws.Range("D$rowNum").Value = $r.D
    }
    $ws.Range("G$rowNum").Value = $r.G

    # Row height matches the wrapped text content, same as other rule rows.
    $ws.Rows.Item($rowNum).RowHeight = $r.Height
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Update the view so the last added row is visible/selected, matching
#    the author's final on-screen state.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("G67").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
